# Apply the edit described by the commit "change template in french":
# rename sheets & the place_of_origin category labels from
# yes/no to always_lived/displaced, and refresh the dependent
# weighted-percentage figures in both subtables.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename the two worksheets.
$ws1.Name = "place_of_origin_always_lived"
$ws2.Name = "place_of_origin_displaced"

# 2) Update the category label shown in column C for every data row
#    on each sheet (place_of_origin: yes -> always_lived, no -> displaced).
for ($r = 4; $r -le 25; $r++) {
    $ws1.Range("C$r").Value = "always_lived"
}
for ($r = 4; $r -le 21; $r++) {
    $ws2.Range("C$r").Value = "displaced"
}

# 3) Refresh the recalculated weighted percentages (columns D:G) that
#    shifted as a result of the regrouping.
$sheet1Changes = @(
    @("D4", 0.4829061080810805),
    @("E4", 0.3086381269252791),
    @("G4", 0.2084557649936404),
    @("D6", 0.6416239440843047),
    @("E6", 0.2285140723300934),
    @("F6", 0.1298619835856019),
    @("D7", 0.2602383032225408),
    @("E7", 0.7397616967774592),
    @("D8", 0.6804654206324127),
    @("E8", 0.1917272943602201),
    @("F8", 0.1069838568764511),
    @("G8", 0.02082342813091612),
    @("D9", 0.7729934045268795),
    @("E9", 0.2270065954731205),
    @("D11", 0.3323660956072528),
    @("E11", 0.4460565073212454),
    @("F11", 0.1107886985357509),
    @("G11", 0.1107886985357509),
    @("D12", 0.487363569267753),
    @("E12", 0.3506459434126348),
    @("F12", 0.1619904873196121),
    @("D13", 0.3073515381473032),
    @("E13", 0.6489498912832892),
    @("F13", 0.04369857056940763),
    @("D14", 0.3563962960946137),
    @("E14", 0.6177580151522395),
    @("F14", 0.02584568875314694),
    @("D16", 0.6243076641910332),
    @("E16", 0.3112937492681226),
    @("F16", 0.06439858654084409),
    @("D18", 0.3382054949988934),
    @("E18", 0.4717568971489077),
    @("F18", 0.1900376078521991),
    @("D20", 0.2954216765382026),
    @("E20", 0.7045783234617974),
    @("D21", 0.4759933064091681),
    @("E21", 0.4524794689665161),
    @("G21", 0.07152722462431582),
    @("D22", 0.1245999338763479),
    @("E22", 0.4108892353518597),
    @("F22", 0.4645108307717923),
    @("D23", 0.4534622879159473),
    @("E23", 0.5465377120840528),
    @("D24", 0.504052662292603),
    @("E24", 0.495947337707397),
    @("D25", 0.3695248772070674),
    @("E25", 0.6304751227929326)
)

foreach ($change in $sheet1Changes) {
    $ws1.Range($change[0]).Value = $change[1]
}

$sheet2Changes = @(
    @("D4", 0.3299638470180129),
    @("E4", 0.6120206706290734),
    @("F4", 0.05801548235291366),
    @("D5", 0.2229487521213348),
    @("E5", 0.6949161349003545),
    @("F5", 0.06217968706701102),
    @("G5", 0.01995542591129953),
    @("D6", 0.06560313118393259),
    @("E6", 0.7259375433774434),
    @("F6", 0.2084593254386241),
    @("D7", 0.3067876531704871),
    @("E7", 0.465825636206937),
    @("F7", 0.2082241078966419),
    @("G7", 0.01916260272593405),
    @("D9", 0.08499592622833563),
    @("E9", 0.7448239379229042),
    @("F9", 0.1701801358487602),
    @("D10", 0.2589672082037577),
    @("E10", 0.3659133674944501),
    @("F10", 0.3751194243017922),
    @("D11", 0.3610955725964223),
    @("E11", 0.6389044274035777),
    @("D12", 0.2022731107742338),
    @("E12", 0.6337914234624509),
    @("F12", 0.1639354657633152),
    @("D15", 0.913650255410942),
    @("E15", 0.08634974458905806),
    @("D18", 0.2231053717126787),
    @("E18", 0.7768946282873211),
    @("D19", 0.2746318732696291),
    @("E19", 0.7253681267303709),
    @("E21", 0.6184791088452267),
    @("F21", 0.1907604455773866),
    @("G21", 0.1907604455773866)
)

foreach ($change in $sheet2Changes) {
    $ws2.Range($change[0]).Value = $change[1]
}

$wb.Save()

